# Re-order the data rows (2-13) of the "Artfynd" sheet according to the
# target revision. Each destination row receives the *entire* content of a
# source row (all populated columns A..AY); the mapping below was derived
# by comparing the before/after cell values for column A (the "Id" column)
# row by row.
#
# destination row -> source row (both referring to row numbers in the
# ORIGINAL/before layout)
#   2  <- 12
#   3  <- 2
#   4  <- 3
#   5  <- 4
#   6  <- 13
#   7  <- 5
#   8  <- 6
#   9  <- 7
#   10 <- 8
#   11 <- 9
#   12 <- 10
#   13 <- 11

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 13
$lastCol = 51   # column AY

$mapping = @{
    2  = 12
    3  = 2
    4  = 3
    5  = 4
    6  = 13
    7  = 5
    8  = 6
    9  = 7
    10 = 8
    11 = 9
    12 = 10
    13 = 11
}

# ---------------------------------------------------------------------
# Phase 1: capture the current contents (value + .NET type) of every
# cell in every source row, BEFORE anything gets overwritten.
# ---------------------------------------------------------------------
$captured = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowData = @{}
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value2
        if ($v -eq $null) {
            $rowData[$c] = $null
        } else {
            $rowData[$c] = @{ Value = $v; TypeName = $v.GetType().FullName }
        }
    }
    $captured[$r] = $rowData
}

# ---------------------------------------------------------------------
# Phase 2: clear the whole block so that columns which are not present
# in the incoming source row (e.g. "M", "K", "L", "N", "AC", "AI") do
# not keep stale leftover values. ClearContents (without touching
# NumberFormat) removes the cells entirely instead of leaving behind
# empty-but-styled cell records.
# ---------------------------------------------------------------------
$clearRange = $ws.Range($ws.Cells.Item($firstRow, 1), $ws.Cells.Item($lastRow, $lastCol))
$clearRange.ClearContents()

# ---------------------------------------------------------------------
# Phase 3: write the captured source-row data into its new destination
# row, restoring the correct cell type as we go (text values that look
# like dates/times must be forced to Text format first so Excel does
# not reinterpret them as date/time serial numbers).
# ---------------------------------------------------------------------
foreach ($destRow in ($mapping.Keys | Sort-Object)) {
    $srcRow = $mapping[$destRow]
    $rowData = $captured[$srcRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $entry = $rowData[$c]
        if ($entry -eq $null) {
            continue
        }
        $cell = $ws.Cells.Item($destRow, $c)
        if ($entry.TypeName -eq "System.String") {
            if ([string]$entry.Value -ne "") {
                $cell.NumberFormat = "@"
                $cell.Value2 = [string]$entry.Value
            }
        } elseif ($entry.TypeName -eq "System.Boolean") {
            $cell.Value2 = [bool]$entry.Value
        } else {
            $cell.Value2 = $entry.Value
        }
    }
}

Write-Host "Row reordering complete."
